$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update report title/date-range text (Volume/Number and week-covering dates) ---
$ws.Range("A8").Value = "Volume 32   Number  48"
$ws.Range("C9").Value = "Report Covering the Week  11/24/2025  Through  11/30/2025"

# --- Update crime statistics table (rows 14-33) with newly collected figures ---
$ws.Range("D14").Copy($ws.Range("C14"))
$ws.Range("C15").Value = 3
$ws.Range("C15").NumberFormat = "#,##0"
$ws.Range("D15").Value = 1
$ws.Range("D15").NumberFormat = "#,##0"
$ws.Range("E15").Value = 200
$ws.Range("E15").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("F15").Value = 4
$ws.Range("H15").Value = 300
$ws.Range("I15").Value = 37
$ws.Range("J15").Value = 32
$ws.Range("K15").Value = 15.625
$ws.Range("L15").Value = 94.736842105263
$ws.Range("M15").Value = 68.181818181818
$ws.Range("N15").Value = -51.315789473684
$ws.Range("C16").Value = 7
$ws.Range("D16").Value = 8
$ws.Range("E16").Value = -12.5
$ws.Range("F16").Value = 29
$ws.Range("G16").Value = 32
$ws.Range("H16").Value = -9.375
$ws.Range("I16").Value = 390
$ws.Range("J16").Value = 366
$ws.Range("K16").Value = 6.557377049180
$ws.Range("L16").Value = 20.743034055727
$ws.Range("M16").Value = -3.703703703703
$ws.Range("N16").Value = -80.071538068472
$ws.Range("C17").Value = 19
$ws.Range("D17").Value = 15
$ws.Range("E17").Value = 26.666666666666
$ws.Range("F17").Value = 79
$ws.Range("G17").Value = 69
$ws.Range("H17").Value = 14.492753623188
$ws.Range("I17").Value = 816
$ws.Range("J17").Value = 826
$ws.Range("K17").Value = -1.210653753026
$ws.Range("L17").Value = 11.323328785811
$ws.Range("M17").Value = 108.163265306122
$ws.Range("N17").Value = -15.088449531737
$ws.Range("C18").Value = 3
$ws.Range("D14").Copy($ws.Range("D18"))
$ws.Range("E14").Copy($ws.Range("E18"))
$ws.Range("F18").Value = 25
$ws.Range("G18").Value = 15
$ws.Range("H18").Value = 66.666666666666
$ws.Range("I18").Value = 192
$ws.Range("K18").Value = -4
$ws.Range("L18").Value = 29.729729729729
$ws.Range("M18").Value = 1.052631578947
$ws.Range("N18").Value = -91.090487238979
$ws.Range("C19").Value = 15
$ws.Range("D19").Value = 18
$ws.Range("E19").Value = -16.666666666666
$ws.Range("F19").Value = 43
$ws.Range("G19").Value = 52
$ws.Range("H19").Value = -17.307692307692
$ws.Range("I19").Value = 604
$ws.Range("J19").Value = 636
$ws.Range("K19").Value = -5.031446540880
$ws.Range("L19").Value = 21.042084168336
$ws.Range("M19").Value = 116.487455197133
$ws.Range("N19").Value = -32.589285714285
$ws.Range("C20").Value = 7
$ws.Range("D20").Value = 5
$ws.Range("E20").Value = 40
$ws.Range("F20").Value = 27
$ws.Range("H20").Value = 68.75
$ws.Range("I20").Value = 240
$ws.Range("J20").Value = 199
$ws.Range("K20").Value = 20.603015075376
$ws.Range("L20").Value = -5.882352941176
$ws.Range("M20").Value = 114.285714285714
$ws.Range("N20").Value = -70.479704797048
$ws.Range("D21").Value = 47
$ws.Range("E21").Value = 14.893617021276
$ws.Range("F21").Value = 209
$ws.Range("G21").Value = 185
$ws.Range("H21").Value = 12.972972972973
$ws.Range("I21").Value = 2287
$ws.Range("J21").Value = 2284
$ws.Range("K21").Value = 0.131348511383
$ws.Range("L21").Value = 14.924623115577
$ws.Range("M21").Value = 61.968838526912
$ws.Range("N21").Value = -66.936533179123
$ws.Range("D14").Copy($ws.Range("D22"))
$ws.Range("E14").Copy($ws.Range("E22"))
$ws.Range("D14").Copy($ws.Range("F22"))
$ws.Range("H22").Value = -100
$ws.Range("M22").Value = -18.75
$ws.Range("C23").Value = 1
$ws.Range("F23").Value = 3
$ws.Range("H23").Value = 50
$ws.Range("I23").Value = 37
$ws.Range("K23").Value = 42.307692307692
$ws.Range("L23").Value = 32.142857142857
$ws.Range("M23").Value = 54.166666666666
$ws.Range("C24").Value = 9
$ws.Range("D24").Value = 15
$ws.Range("E24").Value = -40
$ws.Range("F24").Value = 58
$ws.Range("G24").Value = 74
$ws.Range("H24").Value = -21.621621621621
$ws.Range("I24").Value = 884
$ws.Range("J24").Value = 885
$ws.Range("K24").Value = -0.112994350282
$ws.Range("L24").Value = -7.142857142857
$ws.Range("M24").Value = 3.755868544600
$ws.Range("D25").Value = 7
$ws.Range("E25").Value = -57.142857142857
$ws.Range("G25").Value = 20
$ws.Range("H25").Value = -50
$ws.Range("I25").Value = 227
$ws.Range("J25").Value = 311
$ws.Range("K25").Value = -27.009646302250
$ws.Range("L25").Value = 14.646464646464
$ws.Range("C26").Value = 35
$ws.Range("D26").Value = 25
$ws.Range("E26").Value = 40
$ws.Range("F26").Value = 80
$ws.Range("H26").Value = -2.439024390243
$ws.Range("I26").Value = 969
$ws.Range("J26").Value = 1048
$ws.Range("K26").Value = -7.538167938931
$ws.Range("L26").Value = -11.748633879781
$ws.Range("M26").Value = -18.296795952782
$ws.Range("C27").Value = 3
$ws.Range("C27").NumberFormat = "#,##0"
$ws.Range("D27").Value = 1
$ws.Range("D27").NumberFormat = "#,##0"
$ws.Range("E27").Value = 200
$ws.Range("E27").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("F27").Value = 4
$ws.Range("G27").Value = 1
$ws.Range("H27").Value = 300
$ws.Range("I27").Value = 50
$ws.Range("J27").Value = 46
$ws.Range("K27").Value = 8.695652173913
$ws.Range("L27").Value = 0
$ws.Range("D14").Copy($ws.Range("C28"))
$ws.Range("D28").Value = 2
$ws.Range("E28").Value = -100
$ws.Range("G28").Value = 5
$ws.Range("H28").Value = 0
$ws.Range("J28").Value = 88
$ws.Range("K28").Value = -5.681818181818
$ws.Range("L28").Value = -1.190476190476
$ws.Range("C29").Value = 2
$ws.Range("D14").Copy($ws.Range("D29"))
$ws.Range("E14").Copy($ws.Range("E29"))
$ws.Range("F29").Value = 4
$ws.Range("G29").Value = 5
$ws.Range("H29").Value = -20
$ws.Range("I29").Value = 28
$ws.Range("K29").Value = -56.25
$ws.Range("L29").Value = -47.169811320754
$ws.Range("M29").Value = -46.153846153846
$ws.Range("N29").Value = -83.333333333333
$ws.Range("D14").Copy($ws.Range("D30"))
$ws.Range("E14").Copy($ws.Range("E30"))
$ws.Range("F30").Value = 3
$ws.Range("G30").Value = 3
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 23
$ws.Range("K30").Value = -53.061224489795
$ws.Range("L30").Value = -45.238095238095
$ws.Range("M30").Value = -48.888888888888
$ws.Range("N30").Value = -85.064935064935
$ws.Range("C33").Value = 1
$ws.Range("C33").NumberFormat = "#,##0"
$ws.Range("F33").Value = 1
$ws.Range("F33").NumberFormat = "#,##0"
$ws.Range("I33").Value = 3
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = -40
